$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# --- sheet 1 (sheet1) numeric / text cell updates ---
$ws1.Range('F2').Value = 304
$ws1.Range('F3').Value = 1068
$ws1.Range('F4').Value = 9496
$ws1.Range('F8').Value = 6500
$ws1.Range('F10').Value = 81
$ws1.Range('F11').Value = 9960
$ws1.Range('F12').Value = 11433
$ws1.Range('F13').Value = 1248
$ws1.Range('F14').Value = 1174
$ws1.Range('F15').Value = 4993
$ws1.Range('F16').Value = 811
$ws1.Range('F17').Value = 478
$ws1.Range('F19').Value = 335
$ws1.Range('F22').Value = 271
$ws1.Range('F23').Value = 1881
$ws1.Range('F24').Value = 905
$ws1.Range('F25').Value = 1287
$ws1.Range('F28').Value = 2064
$ws1.Range('F29').Value = 437
$ws1.Range('F30').Value = 647
$ws1.Range('F31').Value = 2712
$ws1.Range('E32').Value = '2024.08.03 10:00-08.04 17:00'
$ws1.Range('F32').Value = 194
$ws1.Range('I32').Value = '//i2.hdslb.com/bfs/openplatform/202407/pSF6KISk1720079182748.jpeg'
$ws1.Range('F33').Value = 1795
$ws1.Range('F35').Value = 814
$ws1.Range('F36').Value = 71
$ws1.Range('F37').Value = 925
$ws1.Range('F38').Value = 43
$ws1.Range('F39').Value = 3381
$ws1.Range('F41').Value = 91
$ws1.Range('F42').Value = 523
$ws1.Range('F43').Value = 591
$ws1.Range('F45').Value = 902
$ws1.Range('F46').Value = 248
$ws1.Range('F48').Value = 4223
$ws1.Range('F49').Value = 65

# --- sheet 2 (sheet2) updates ---
$ws2.Range('F9').Value = 21

# Insert a new row at position 29 (演出 sheet), shifting the existing
# row 29 (2024-12-24 event) down to row 30, then populate the new row 29
# with the inserted 2024-12-22 piano concert event.
$ws2.Rows(29).Insert()

$ws2.Range('A29').Value = 28
$ws2.Range('A30').Value = 29
$ws2.Range('B29').NumberFormat = '@'
$ws2.Range('B29').Value = '2024-12-22'
$ws2.Range('C29').Value = '杭州·《你的名字》《天气之子》《铃芽之旅》——新海诚动漫三部曲钢琴演奏会'
$ws2.Range('D29').Value = '曙光路31号 浙江音乐厅'
$ws2.Range('E29').Value = '2024.12.22 19:30-12.22 21:00'
$ws2.Range('F29').Value = 0
$ws2.Range('G29').Value = 153
$ws2.Range('H29').Value = 'https://show.bilibili.com/platform/detail.html?id=88648'
$ws2.Range('I29').Value = '//i1.hdslb.com/bfs/openplatform/202407/nEB3TPxP1720064877363.jpeg'

# --- sheet 3 (sheet3) numeric / text cell updates ---
$ws3.Range('F2').Value = 6014

# --- sheet 4 (sheet4) numeric / text cell updates ---
$ws4.Range('F2').Value = 304
$ws4.Range('F3').Value = 1068
$ws4.Range('F4').Value = 9496
$ws4.Range('F10').Value = 6500
$ws4.Range('F11').Value = 9960
$ws4.Range('F12').Value = 11434
$ws4.Range('F13').Value = 1174
$ws4.Range('F14').Value = 811
$ws4.Range('F15').Value = 478
$ws4.Range('F17').Value = 335
$ws4.Range('F18').Value = 21
$ws4.Range('F23').Value = 271
$ws4.Range('F24').Value = 1881
$ws4.Range('F25').Value = 905
$ws4.Range('F26').Value = 1287
$ws4.Range('F28').Value = 2064
$ws4.Range('F29').Value = 647
$ws4.Range('F30').Value = 2712
$ws4.Range('E31').Value = '2024.08.03 10:00-08.04 17:00'
$ws4.Range('F31').Value = 194
$ws4.Range('I31').Value = '//i2.hdslb.com/bfs/openplatform/202407/pSF6KISk1720079182748.jpeg'
$ws4.Range('F33').Value = 814
$ws4.Range('F38').Value = 71
$ws4.Range('F40').Value = 43
$ws4.Range('F42').Value = 91
$ws4.Range('F43').Value = 523
$ws4.Range('F44').Value = 591
$ws4.Range('F45').Value = 902
$ws4.Range('F46').Value = 248
$ws4.Range('F48').Value = 4223

